$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = 'FULL'
$ws.Range("C2").Value = ''
$ws.Range("D2").Value = 'JohnFull'
$ws.Range("E2").Value = 'DoeFull'
$ws.Range("F2").Value = '1111111111'
$ws.Range("G2").Value = 'johndoefull@gmail.com'
$ws.Range("H2").Value = '1901-01-01T01:01:01.111Z'
$ws.Range("I2").Value = 'Education 1'
$ws.Range("J2").Value = 'HospitalName 1'
$ws.Range("K2").Value = 'HospitalStreet 1'
$ws.Range("L2").Value = '111111'
$ws.Range("M2").Value = 'HospitalCity 1'

$ws.Range("A3").Value = 'MANDATORY'
$ws.Range("C3").Value = ''
$ws.Range("D3").Value = 'JohnMandatory'
$ws.Range("E3").Value = 'DoeMandatory'
$ws.Range("F3").Value = '2222222222'
$ws.Range("G3").Value = 'johndoemandatory@gmail.com'
$ws.Range("H3").Value = '1902-02-02T02:02:02.222Z'
$ws.Range("I3").Value = 'Education 2'
$ws.Range("J3").Value = 'HospitalName 2'
$ws.Range("K3").Value = 'HospitalStreet 2'
$ws.Range("L3").Value = '222222'
$ws.Range("M3").Value = 'HospitalCity 2'

$ws.Range("A4").Value = 'ADDITIONAL'
$ws.Range("B4").Value = '3333333333'
$ws.Range("C4").Value = '3333333333'
$ws.Range("D4").Value = 'JohnAdditional'
$ws.Range("E4").Value = 'DoeAdditional'
$ws.Range("F4").Value = '3333333333'
$ws.Range("G4").Value = 'johndoeadditional@gmail.com'
$ws.Range("H4").Value = '1903-03-03T03:03:03.333Z'
$ws.Range("I4").Value = 'Education 3'
$ws.Range("J4").Value = 'HospitalName 3'
$ws.Range("K4").Value = 'HospitalStreet 3'
$ws.Range("L4").Value = '333333'
$ws.Range("M4").Value = 'HospitalCity 3'

$ws.Range("A5").Value = 'INVALID'
$ws.Range("C5").Value = ''
$ws.Range("D5").Value = 'JohnInvalidPincode'
$ws.Range("E5").Value = 'DoeInvalidPincode'
$ws.Range("F5").Value = '4444444444'
$ws.Range("G5").Value = 'johndoeinvalid@gmail.com'
$ws.Range("H5").Value = '1904-04-04T04:04:04.444Z'
$ws.Range("I5").Value = 'Education 4'
$ws.Range("J5").Value = 'HospitalName 4'
$ws.Range("K5").Value = 'HospitalStreet 4'
$ws.Range("L5").Value = '444444444444444444'
$ws.Range("M5").Value = 'HospitalCity 4'

$ws.Range("A6").Value = 'GET_BY_ID'
$ws.Range("C6").Value = ''
$ws.Range("D6").Value = 'JohnGetById'
$ws.Range("E6").Value = 'DoeGetById'
$ws.Range("F6").Value = '5555555555'
$ws.Range("G6").Value = 'johndoegetbyid@gmail.com'
$ws.Range("H6").Value = '1905-05-05T05:05:05.555Z'
$ws.Range("I6").Value = 'Education 5'
$ws.Range("J6").Value = 'HospitalName 5'
$ws.Range("K6").Value = 'HospitalStreet 5'
$ws.Range("L6").Value = '555555'
$ws.Range("M6").Value = 'HospitalCity 5'

$ws.Range("A7").Value = 'GET_ALL'
$ws.Range("C7").Value = ''
$ws.Range("D7").Value = 'JohnGetAll'
$ws.Range("E7").Value = 'DoeGetAll'
$ws.Range("F7").Value = '6666666666'
$ws.Range("G7").Value = 'johndoegetall@gmail.com'
$ws.Range("H7").Value = '1906-06-06T06:06:06.666Z'
$ws.Range("I7").Value = 'Education 6'
$ws.Range("J7").Value = 'HospitalName 6'
$ws.Range("K7").Value = 'HospitalStreet 6'
$ws.Range("L7").Value = '666666'
$ws.Range("M7").Value = 'HospitalCity 6'

$ws.Range("A8").Value = 'DELETE_BY_ID'
$ws.Range("C8").Value = ''
$ws.Range("D8").Value = 'JohnDeleteById'
$ws.Range("E8").Value = 'DoeDeleteById'
$ws.Range("F8").Value = '7777777777'
$ws.Range("G8").Value = 'johndoedeletebyid@gmail.com'
$ws.Range("H8").Value = '1907-07-07T07:07:07.777Z'
$ws.Range("I8").Value = 'Education 7'
$ws.Range("J8").Value = 'HospitalName 7'
$ws.Range("K8").Value = 'HospitalStreet 7'
$ws.Range("L8").Value = '777777'
$ws.Range("M8").Value = 'HospitalCity 7'

